# Update gh-pages to output generated at 456a3b4
# Applies numeric "want to go" (column F) count changes across sheets.

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3,6).Value = 315
$ws1.Cells.Item(4,6).Value = 1293
$ws1.Cells.Item(6,6).Value = 332
$ws1.Cells.Item(7,6).Value = 1136
$ws1.Cells.Item(8,6).Value = 438
$ws1.Cells.Item(9,6).Value = 7018
$ws1.Cells.Item(13,6).Value = 7905
$ws1.Cells.Item(15,6).Value = 50
$ws1.Cells.Item(16,6).Value = 5472
$ws1.Cells.Item(17,6).Value = 47
$ws1.Cells.Item(18,6).Value = 2364
$ws1.Cells.Item(20,6).Value = 4546
$ws1.Cells.Item(21,6).Value = 286
$ws1.Cells.Item(23,6).Value = 77
$ws1.Cells.Item(25,6).Value = 345
$ws1.Cells.Item(26,6).Value = 243
$ws1.Cells.Item(27,6).Value = 8
$ws1.Cells.Item(28,6).Value = 2206
$ws1.Cells.Item(29,6).Value = 22
$ws1.Cells.Item(30,6).Value = 253
$ws1.Cells.Item(32,6).Value = 96
$ws1.Cells.Item(33,6).Value = 560
$ws1.Cells.Item(36,6).Value = 1450
$ws1.Cells.Item(39,6).Value = 2219
$ws1.Cells.Item(40,6).Value = 2200
$ws1.Cells.Item(41,6).Value = 1
$ws1.Cells.Item(42,6).Value = 1

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3,6).Value = 70
$ws2.Cells.Item(4,6).Value = 51

# 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3,6).Value = 1270

# 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4,6).Value = 1270
$ws4.Cells.Item(5,6).Value = 315
$ws4.Cells.Item(7,6).Value = 332
$ws4.Cells.Item(8,6).Value = 1136
$ws4.Cells.Item(9,6).Value = 438
$ws4.Cells.Item(10,6).Value = 7018
$ws4.Cells.Item(14,6).Value = 7905
$ws4.Cells.Item(16,6).Value = 50
$ws4.Cells.Item(17,6).Value = 5472
$ws4.Cells.Item(18,6).Value = 47
$ws4.Cells.Item(19,6).Value = 2364
$ws4.Cells.Item(21,6).Value = 4546
$ws4.Cells.Item(22,6).Value = 286
$ws4.Cells.Item(24,6).Value = 77
$ws4.Cells.Item(25,6).Value = 70
$ws4.Cells.Item(27,6).Value = 51
$ws4.Cells.Item(28,6).Value = 345
$ws4.Cells.Item(29,6).Value = 243
$ws4.Cells.Item(30,6).Value = 8
$ws4.Cells.Item(31,6).Value = 2206
$ws4.Cells.Item(32,6).Value = 22
$ws4.Cells.Item(33,6).Value = 253
$ws4.Cells.Item(35,6).Value = 96
$ws4.Cells.Item(36,6).Value = 560
$ws4.Cells.Item(40,6).Value = 1450
$ws4.Cells.Item(43,6).Value = 2219
$ws4.Cells.Item(45,6).Value = 2200
$ws4.Cells.Item(46,6).Value = 1
$ws4.Cells.Item(47,6).Value = 1
